$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently holds line1..line6 (rows 2-7) followed by extr1..extr8
# (rows 8-15). Two new lines (line7, line8) need to be inserted right after
# line6, so the extr1..extr8 block has to move down by two rows (from rows
# 8-15 to rows 10-17). Shift the existing data down cell-by-cell (bottom to
# top so nothing gets clobbered) instead of using a structural row Insert,
# which would otherwise introduce an unwanted extra cell style.
for ($r = 15; $r -ge 8; $r--) {
    for ($c = 1; $c -le 5; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 2, $c).Value = $v
    }
}

# Rows 16 and 17 are brand new, so they don't yet carry the bold/centered/
# bordered formatting used for column A throughout the table. Copy that
# formatting from an existing row (reuses the existing style, adds no new
# style entries).
$ws.Cells.Item(7, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Write the new line7/line8 rows (8-9) and refresh the name labels for the
# (now shifted) extr1..extr8 rows (10-17).
$names = @("line7", "line8", "extr1", "extr2", "extr3", "extr4", "extr5", "extr6", "extr7", "extr8")
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $names[$i]
}

# Updated from_bus / to_bus / in_service results for rows 8-17.
$ws.Cells.Item(8, 3).Value  = 14
$ws.Cells.Item(8, 4).Value  = 11
$ws.Cells.Item(8, 5).Value  = $true

$ws.Cells.Item(9, 3).Value  = 16
$ws.Cells.Item(9, 4).Value  = 9
$ws.Cells.Item(9, 5).Value  = $false

$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
